$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) GTRI / MIDAS bullet: the old description was three runs that spelled out
#    "(MIDAS) Machine Intelligence-Directed Attack Simulator." — replace the
#    whole sentence with the new framing. The new wording/formatting exactly
#    matches another occurrence of this same sentence that already exists
#    further down in the resume (another job's bullet), so we copy that
#    run's FormattedText over the old runs to pick up its exact character
#    formatting (drops the stray w:bCs, picks up w:lang eastAsia="zh-TW").
# ---------------------------------------------------------------------------
$midasPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*Machine Intelligence-Directed Attack Simulator*") {
        $midasPara = $cand
        break
    }
}

$srcRange = $d.Content
$srcRange.Find.Execute("Adversarial Machine Learning for PDF Malware.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$srcFormatted = $srcRange.FormattedText

$dstRange = $midasPara.Range
$dstRange.Find.Execute("(MIDAS) Machine Intelligence-Directed Attack Simulator.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$dstRange.FormattedText = $srcFormatted

# ---------------------------------------------------------------------------
# 2) "Add Integrated Gradients attribution method to MIDAS, a platform for
#    studying AI security." -> "Implement the Integrated Gradients
#    attribution method."
# ---------------------------------------------------------------------------
$igParaIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*Integrated Gradients*") {
        $igParaIndex = $i
        break
    }
}
$igPara = $d.Paragraphs.Item($igParaIndex)
$igRange = $igPara.Range
$igRange.Find.Execute("Add Integrated Gradients attribution method to MIDAS, a platform for studying AI security.", $true, $false, $false, $false, $false, $true, 1, $false, "Implement the Integrated Gradients attribution method.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Selecting features for PDF malware detectors to ensure robustness
#    against adversarial attacks." -> "Select features with robustness
#    against adversarial attacks for PDF malware detectors."
# ---------------------------------------------------------------------------
$featParaIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*Selecting features*") {
        $featParaIndex = $i
        break
    }
}
$featPara = $d.Paragraphs.Item($featParaIndex)
$featRange = $featPara.Range
$featRange.Find.Execute("Selecting features for PDF malware detectors to ensure robustness against adversarial attacks.", $true, $false, $false, $false, $false, $true, 1, $false, "Select features with robustness against adversarial attacks for PDF malware detectors.", 2) | Out-Null
